# "chot csdl lan 2" — add start/end date columns to the lop_hoc sheet and
# switch the trangthai values from raw codes to their Vietnamese display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lop_hoc")

# Insert two new columns (ngaybatdau / ngayketthuc) right before the
# existing khoa_hoc_id column (column C), pushing khoa_hoc_id/trangthai/
# ghichu out to E/F/G.
$ws.Range("C1:D1").EntireColumn.Insert()

# The new columns hold date-like text ("2020/05/07"), so format them as
# Text up front (reuses the workbook's existing Text style) before writing
# any values -- otherwise Excel would auto-coerce the strings into date
# serial numbers.
$ws.Range("C1:D3").NumberFormat = "@"

# trangthai: replace the old status codes with their Vietnamese labels.
# (written first so the new shared-string entries land in the same order
# the workbook ends up with)
$ws.Range("F2").Value = "Đang hoạt động"
$ws.Range("F3").Value = "Đã kết thúc"

# ngaybatdau / ngayketthuc data
$ws.Range("C2").Value = "2020/05/07"
$ws.Range("D2").Value = "2020/07/07"
$ws.Range("C3").Value = "2020/06/08"
$ws.Range("D3").Value = "2020/08/15"

# Headers for the two new columns
$ws.Range("C1").Value = "ngaybatdau"
$ws.Range("D1").Value = "ngayketthuc"

# Match the column widths used elsewhere in the sheet as closely as possible
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 16.140625

# Leave the selection on the new ngayketthuc header, matching the saved view
$ws.Range("D1").Select()
